$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# FuncLoc (AV2) and Previous Doc (AX2) are already Text-formatted cells,
# so a plain value assignment keeps them as shared-string text cells.
$ws.Range("AV2").Value = "ABCD391135"
$ws.Range("AX2").Value = "4039607280"

# SAID (AW2) is Number-formatted (s=16, numFmtId "0"); the target keeps
# that same style but stores the value as text, so we momentarily force
# Text formatting to write it as a string, then restore the original
# number format code so the cell's style/index is preserved.
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = "2425289268"
$ws.Range("AW2").NumberFormat = "0"
